$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: update title and link
$ws.Range("D28").Value = "WSL2 ROS (with Docker)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/161"

# Row 37: update title and link
$ws.Range("D37").Value = "[Paper Review] RepSum: Unsupervised Dialogue Summarization based on Replacement Strategy"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1825&mod=document&pageid=1"
